# Weekly crime-stat refresh: shift the reporting week forward by one week
# (Volume/Number label + date range) and replace the precinct crime figures
# with the newly collected week's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (rich-text shared strings) ---------------------------
# "Volume 30   Number  19" -> "...  20"
$ws.Range("A8").Characters(21, 2).Text = "20"

# "Report Covering the Week  5/8/2023  Through  5/14/2023"
#   -> "...5/15/2023  Through  5/21/2023"
# Replace the later substring first so the earlier offset stays valid.
$ws.Range("C9").Characters(46, 9).Text = "5/21/2023"
$ws.Range("C9").Characters(27, 8).Text = "5/15/2023"

# --- Cells that change between numeric and text ("0" / "***.*") forms -------
# Copying from a cell already in the desired shape carries over both the
# shared-string value and the style index, matching how Excel itself would
# flip a cell between "empty-looking" text and a real number.
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E15").Copy($ws.Range("H22"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("C29"))

$ws.Range("C26").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 3

# --- Weekly figures refresh (Week-to-Date / 28-Day / YTD / comparison %s) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -78.571428571428

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -5.882352941176
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = 3.076923076923
$ws.Range("L16").Value = 9.836065573770
$ws.Range("M16").Value = -26.373626373626
$ws.Range("N16").Value = -77.516778523489

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 33.333333333333
$ws.Range("I17").Value = 136
$ws.Range("J17").Value = 122
$ws.Range("K17").Value = 11.475409836065
$ws.Range("L17").Value = 12.396694214876
$ws.Range("M17").Value = 94.285714285714
$ws.Range("N17").Value = -35.849056603773

$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -31.25
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = -27.692307692307
$ws.Range("L18").Value = 51.612903225806
$ws.Range("M18").Value = 104.347826086957
$ws.Range("N18").Value = -69.677419354838

$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -53.846153846153
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -29.411764705882
$ws.Range("I19").Value = 156
$ws.Range("J19").Value = 157
$ws.Range("K19").Value = -0.636942675159
$ws.Range("L19").Value = 31.092436974789
$ws.Range("M19").Value = 79.310344827586
$ws.Range("N19").Value = -33.333333333333

$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -57.142857142857
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 123.529411764706
$ws.Range("N20").Value = -73.611111111111

$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -37.142857142857
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -23.931623931623
$ws.Range("I21").Value = 455
$ws.Range("J21").Value = 459
$ws.Range("K21").Value = -0.871459694989
$ws.Range("L21").Value = 26.388888888888
$ws.Range("M21").Value = 54.761904761904
$ws.Range("N21").Value = -58.103130755064

$ws.Range("C23").Value = 9
$ws.Range("E23").Value = -18.181818181818
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 140
$ws.Range("J23").Value = 146
$ws.Range("K23").Value = -4.109589041095
$ws.Range("L23").Value = 0.719424460431
$ws.Range("M23").Value = 68.674698795180

$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -5.263157894736
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 22.058823529411
$ws.Range("I24").Value = 346
$ws.Range("J24").Value = 283
$ws.Range("K24").Value = 22.261484098939
$ws.Range("L24").Value = 16.107382550335
$ws.Range("M24").Value = 45.991561181434

$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = -30.357142857142
$ws.Range("I25").Value = 213
$ws.Range("J25").Value = 222
$ws.Range("K25").Value = -4.054054054054
$ws.Range("L25").Value = 6.5
$ws.Range("M25").Value = -17.760617760617

$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 10
$ws.Range("K26").Value = -9.090909090909
$ws.Range("L26").Value = 11.111111111111

$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = 9.523809523809
$ws.Range("L27").Value = -17.857142857142

$ws.Range("L28").Value = -33.333333333333

$ws.Range("L29").Value = -36.363636363636
